# daily auto push: 2026-01-18 18:40 UTC
# Insert two new rows of data at row 651 (pushing the existing rows down by
# two, e.g. old row 651 -> new row 653, ..., old row 692 -> new row 694)
# and populate the two new rows with the day's readings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows above the current row 651, shifting everything
# from old-651..old-692 down to new-653..new-694.
$ws.Rows.Item(651).Resize(2).Insert()

# New row 651: 2026/01/18, 日, 22, 18
$ws.Cells.Item(651, 1).Value2 = "'2026/01/18"
$ws.Cells.Item(651, 1).ClearFormats()
$ws.Cells.Item(651, 2).Value2 = "日"
$ws.Cells.Item(651, 3).Value2 = 22
$ws.Cells.Item(651, 4).Value2 = 18

# New row 652: 2026/01/19, 月, 1, 19
$ws.Cells.Item(652, 1).Value2 = "'2026/01/19"
$ws.Cells.Item(652, 1).ClearFormats()
$ws.Cells.Item(652, 2).Value2 = "月"
$ws.Cells.Item(652, 3).Value2 = 1
$ws.Cells.Item(652, 4).Value2 = 19
